# ✅ complete Level2Scene features
#
# Highlights the QA checklist bullets on the Level1Scene, Level2Scene, and
# DemoScene "Step" slides (mostly green, one yellow), and merges the split
# "DemoScene." / "c" runs on the DemoScene slide into a single run.

function Set-ParaHighlights {
    param($TextRange, $ColorsByIndex)

    $full = $TextRange.Text
    $parts = $full -split "`r"
    $pos = 1
    for ($i = 0; $i -lt $parts.Count; $i++) {
        $part = $parts[$i]
        $len = $part.Length
        if ($len -gt 0 -and $ColorsByIndex.ContainsKey($i)) {
            $color = $ColorsByIndex[$i]
            $sub = $TextRange.Characters($pos, $len)
            $sub.Font.Highlight = $color
        }
        $pos += $len + 1
    }
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 13 - "Step 10" / Level1Scene.c
# ---------------------------------------------------------------------
$slide13 = $p.Slides.Item(13)
$tr13 = $slide13.Shapes.Item("Content Placeholder 2").TextFrame.TextRange
$colors13 = @{
    1 = "00FF00"
    2 = "00FF00"
    3 = "FFFF00"
    4 = "00FF00"
    5 = "00FF00"
    6 = "00FF00"
    7 = "00FF00"
    8 = "00FF00"
}
Set-ParaHighlights $tr13 $colors13

# ---------------------------------------------------------------------
# Slide 14 - "Step 11" / Level2Scene.c
# ---------------------------------------------------------------------
$slide14 = $p.Slides.Item(14)
$tr14 = $slide14.Shapes.Item("Content Placeholder 2").TextFrame.TextRange
$colors14 = @{
    1 = "00FF00"
    2 = "00FF00"
    3 = "00FF00"
    4 = "00FF00"
    5 = "00FF00"
    6 = "00FF00"
    7 = "00FF00"
    8 = "00FF00"
}
Set-ParaHighlights $tr14 $colors14

# ---------------------------------------------------------------------
# Slide 15 - "Step 12" / DemoScene.c
# ---------------------------------------------------------------------
$slide15 = $p.Slides.Item(15)
$tr15 = $slide15.Shapes.Item("Content Placeholder 2").TextFrame.TextRange

# Merge the "DemoScene." run and the "c" run into a single "DemoScene.c"
# run, keeping the second run's properties (dirty="0" err="1").
$firstRunRange = $tr15.Characters(1, 10)
$firstRunRange.Text = ""
$cRunRange = $tr15.Characters(1, 1)
[void]$cRunRange.InsertBefore("DemoScene.")

$colors15 = @{
    1 = "00FF00"
    2 = "00FF00"
    3 = "00FF00"
    4 = "00FF00"
    5 = "00FF00"
    6 = "00FF00"
    7 = "00FF00"
}
Set-ParaHighlights $tr15 $colors15
